# Applies the numeric cell updates from the Mandragora_Profits commit diff
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1464.1459
$ws.Range("J17").Value = 1464.1459
$ws.Range("L17").Value = 4392.4377
$ws.Range("N17").Value = -4728.4377

$ws.Range("H46").Value = 55556892
$ws.Range("J46").Value = 2000
$ws.Range("L46").Value = 6000
$ws.Range("N46").Value = -6238

$ws.Range("H55").Value = 383.2
$ws.Range("I55").Value = 321.93332
$ws.Range("J55").Value = 475.1
$ws.Range("K55").Value = 321.93332
$ws.Range("L55").Value = 475.1
$ws.Range("M55").Value = -107.93332
$ws.Range("N55").Value = -903.1

$ws.Range("H58").Value = 1436.5
$ws.Range("I58").Value = 394.16666
$ws.Range("J58").Value = 3000
$ws.Range("K58").Value = 1182.49998
$ws.Range("L58").Value = 9000
$ws.Range("M58").Value = -1032.49998
$ws.Range("N58").Value = -9300

$ws.Range("H60").Value = 55556892
$ws.Range("J60").Value = 2000
$ws.Range("L60").Value = 6000
$ws.Range("N60").Value = -6968

$ws.Range("H61").Value = 575
$ws.Range("I61").Value = 150
$ws.Range("K61").Value = 450
$ws.Range("M61").Value = -278

$ws.Range("H111").Value = 1391.5
$ws.Range("I111").Value = 1050
$ws.Range("J111").Value = 1733
$ws.Range("K111").Value = 3150
$ws.Range("L111").Value = 5199
$ws.Range("M111").Value = -83
$ws.Range("N111").Value = -11333

$ws.Range("H135").Value = 735.78845
$ws.Range("I135").Value = 273.13333
$ws.Range("J135").Value = 1366.6818
$ws.Range("K135").Value = 2458.19997
$ws.Range("L135").Value = 12300.1362
$ws.Range("M135").Value = 76.80002999999988
$ws.Range("N135").Value = -17370.1362

$ws.Range("H137").Value = 1521.0862
$ws.Range("I137").Value = 2013.8334
$ws.Range("J137").Value = 993.1429000000001
$ws.Range("K137").Value = 6041.5002
$ws.Range("L137").Value = 2979.4287
$ws.Range("M137").Value = -3491.5002
$ws.Range("N137").Value = -8079.4287

$ws.Range("H138").Value = 1172.7174
$ws.Range("I138").Value = 1045.5
$ws.Range("J138").Value = 1533.1666
$ws.Range("K138").Value = 3136.5
$ws.Range("L138").Value = 4599.4998
$ws.Range("M138").Value = 2003.5
$ws.Range("N138").Value = -14879.4998

$ws.Range("H141").Value = 2588.762
$ws.Range("I141").Value = 952.5585
$ws.Range("J141").Value = 20587
$ws.Range("K141").Value = 2857.6755
$ws.Range("L141").Value = 61761
$ws.Range("M141").Value = 2322.3245
$ws.Range("N141").Value = -72121

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1710.7826
$ws.Range("I61").Value = 1782.7
$ws.Range("J61").Value = 1655.4615
$ws.Range("K61").Value = 1782.7
$ws.Range("L61").Value = 1655.4615
$ws.Range("M61").Value = -1570.7
$ws.Range("N61").Value = -2079.4615

$ws.Range("H132").Value = 744.66
$ws.Range("I132").Value = 487.77194
$ws.Range("J132").Value = 1085.186
$ws.Range("K132").Value = 1463.31582
$ws.Range("L132").Value = 3255.558
$ws.Range("M132").Value = 1066.68418
$ws.Range("N132").Value = -8315.558000000001

$ws.Range("H136").Value = 1710.7826
$ws.Range("I136").Value = 1782.7
$ws.Range("J136").Value = 1655.4615
$ws.Range("K136").Value = 5348.1
$ws.Range("L136").Value = 4966.3845
$ws.Range("M136").Value = -2798.1
$ws.Range("N136").Value = -10066.3845

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 434.8889
$ws.Range("I80").Value = 309.7647
$ws.Range("J80").Value = 647.6
$ws.Range("K80").Value = 309.7647
$ws.Range("L80").Value = 647.6
$ws.Range("M80").Value = 688.2353000000001
$ws.Range("N80").Value = -2643.6

$ws.Range("H83").Value = 434.8889
$ws.Range("I83").Value = 309.7647
$ws.Range("J83").Value = 647.6
$ws.Range("K83").Value = 1548.8235
$ws.Range("L83").Value = 3238
$ws.Range("M83").Value = 3443.1765
$ws.Range("N83").Value = -13222

$ws.Range("H134").Value = 4071.7048
$ws.Range("I134").Value = 1490.973
$ws.Range("J134").Value = 8050.3335
$ws.Range("K134").Value = 4472.919
$ws.Range("L134").Value = 24151.0005
$ws.Range("M134").Value = -1937.919
$ws.Range("N134").Value = -29221.0005

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4631492.5
$ws.Range("I31").Value = 1319.9445
$ws.Range("J31").Value = 18522010
$ws.Range("K31").Value = 1319.9445
$ws.Range("L31").Value = 18522010
$ws.Range("M31").Value = -1024.9445
$ws.Range("N31").Value = -18522600

$ws.Range("H34").Value = 4631492.5
$ws.Range("I34").Value = 1319.9445
$ws.Range("J34").Value = 18522010
$ws.Range("K34").Value = 1319.9445
$ws.Range("L34").Value = 18522010
$ws.Range("M34").Value = -1117.9445
$ws.Range("N34").Value = -18522414

$ws.Range("H58").Value = 1003828.44
$ws.Range("I58").Value = 2377.5417
$ws.Range("J58").Value = 1928244.8
$ws.Range("K58").Value = 2377.5417
$ws.Range("L58").Value = 1928244.8
$ws.Range("M58").Value = -2174.5417
$ws.Range("N58").Value = -1928650.8

$ws.Range("H105").Value = 1010.1316
$ws.Range("I105").Value = 1007.4054
$ws.Range("J105").Value = 1111
$ws.Range("K105").Value = 1007.4054
$ws.Range("L105").Value = 1111
$ws.Range("M105").Value = 739.5946
$ws.Range("N105").Value = -4605

$ws.Range("H132").Value = 1740.1571
$ws.Range("I132").Value = 960.8684
$ws.Range("J132").Value = 2665.5625
$ws.Range("K132").Value = 2882.6052
$ws.Range("L132").Value = 7996.6875
$ws.Range("M132").Value = -352.6052
$ws.Range("N132").Value = -13056.6875

$ws.Range("H134").Value = 1086.7604
$ws.Range("I134").Value = 563.7077
$ws.Range("J134").Value = 2183.484
$ws.Range("K134").Value = 1691.1231
$ws.Range("L134").Value = 6550.451999999999
$ws.Range("M134").Value = 843.8768999999998
$ws.Range("N134").Value = -11620.452

$ws.Range("H136").Value = 1003828.44
$ws.Range("I136").Value = 2377.5417
$ws.Range("J136").Value = 1928244.8
$ws.Range("K136").Value = 7132.625100000001
$ws.Range("L136").Value = 5784734.4
$ws.Range("M136").Value = -4582.625100000001
$ws.Range("N136").Value = -5789834.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H60").Value = 1507.8334
$ws.Range("I60").Value = 318.8
$ws.Range("J60").Value = 2357.1428
$ws.Range("K60").Value = 956.4000000000001
$ws.Range("L60").Value = 7071.428400000001
$ws.Range("M60").Value = -705.4000000000001
$ws.Range("N60").Value = -7573.428400000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3104.2104
$ws.Range("I80").Value = 2528.5
$ws.Range("J80").Value = 3743.889
$ws.Range("K80").Value = 2528.5
$ws.Range("L80").Value = 3743.889
$ws.Range("M80").Value = -1530.5
$ws.Range("N80").Value = -5739.889

$ws.Range("H83").Value = 3104.2104
$ws.Range("I83").Value = 2528.5
$ws.Range("J83").Value = 3743.889
$ws.Range("K83").Value = 12642.5
$ws.Range("L83").Value = 18719.445
$ws.Range("M83").Value = -7650.5
$ws.Range("N83").Value = -28703.445

$ws.Range("H132").Value = 652378.4
$ws.Range("I132").Value = 992931.8
$ws.Range("J132").Value = 2231
$ws.Range("K132").Value = 2978795.4
$ws.Range("L132").Value = 6693
$ws.Range("M132").Value = -2976265.4
$ws.Range("N132").Value = -11753

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2399.0513
$ws.Range("I40").Value = 1969.4667
$ws.Range("J40").Value = 3831
$ws.Range("K40").Value = 1969.4667
$ws.Range("L40").Value = 3831
$ws.Range("M40").Value = -1833.4667
$ws.Range("N40").Value = -4103

$ws.Range("H120").Value = 40698
$ws.Range("J120").Value = 40698
$ws.Range("L120").Value = 40698
$ws.Range("N120").Value = -50374

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1193.6666
$ws.Range("I132").Value = 789.9808
$ws.Range("J132").Value = 1917.5172
$ws.Range("K132").Value = 2369.9424
$ws.Range("L132").Value = 5752.5516
$ws.Range("M132").Value = 160.0576000000001
$ws.Range("N132").Value = -10812.5516
